$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the table, A1:E51) used to stage values that look like plain
# numbers (e.g. "583.66") as literal text via a ="..." formula + copy/paste-values,
# so the target cell keeps storing a text string (matching the sheets inlineStr
# cells) instead of Excel auto-converting it to a Number when assigned directly,
# and without leaving any NumberFormat / style residue behind.
$helper = $ws.Range("G1")

$ws.Range("D2").Value = "60.801.19"
$ws.Range("E2").Value = "  -3.38%  "
$ws.Range("D3").Value = "2.914.46"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  -0.07%  "
$helper.Formula = "=""583.66"""
$helper.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -1.73%  "
$helper.Formula = "=""144.87"""
$helper.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -5.79%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").Value = "2.912.60"
$ws.Range("E9").Value = "  -3.83%  "
$helper.Formula = "=""6.84"""
$helper.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("E11").Value = "  -4.83%  "
$ws.Range("E12").Value = "  -4.20%  "
$helper.Formula = "=""0.0000226"""
$helper.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -3.82%  "
$helper.Formula = "=""33.53"""
$helper.Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -6.29%  "
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "3.397.08"
$ws.Range("D17").Value = "60.750.34"
$ws.Range("E17").Value = "  -3.39%  "
$helper.Formula = "=""6.74"""
$helper.Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -4.93%  "
$ws.Range("D19").Value = "2.914.19"
$ws.Range("E19").Value = "  -3.99%  "
$helper.Formula = "=""431.01"""
$helper.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -5.24%  "
$ws.Range("E21").Value = "  -4.85%  "
$helper.Formula = "=""0.682"""
$helper.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("E23").Value = "  -4.85%  "
$helper.Formula = "=""80.25"""
$helper.Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -3.63%  "
$helper.Formula = "=""10.89"""
$helper.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -2.65%  "
$helper.Formula = "=""2.21"""
$helper.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -4.70%  "
$helper.Formula = "=""11.87"""
$helper.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -4.15%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -0.02%  "
$helper.Formula = "=""7.19"""
$helper.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("E31").Value = "  -3.36%  "
$helper.Formula = "=""2.16"""
$helper.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("E34").Value = "  -3.57%  "
$ws.Range("D35").Value = "0.0₃0871"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("E37").Value = "  -4.95%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$helper.Formula = "=""3.02"""
$helper.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -5.61%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$helper.Formula = "=""0.128"""
$helper.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -0.11%  "
$helper.Formula = "=""49.62"""
$helper.Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("E41").Value = "  -5.60%  "
$ws.Range("E42").Value = "  -5.12%  "
$helper.Formula = "=""0.292"""
$helper.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -3.77%  "
$helper.Formula = "=""41.03"""
$helper.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -2.31%  "
$helper.Formula = "=""376.43"""
$helper.Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -5.47%  "
$helper.Formula = "=""0.0349"""
$helper.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -3.30%  "
$ws.Range("D47").Value = "2.671.31"
$ws.Range("E47").Value = "  -2.06%  "
$helper.Formula = "=""132.23"""
$helper.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +0.21%  "
$helper.Formula = "=""24.33"""
$helper.Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -2.10%  "

$helper.Value = ""
